# This commit ("Moving from 2.0.1 to 2.0.2") re-saved the template with a
# newer tool version. Diffing the canonical OOXML shows every changed line
# is a pure XML-attribute/namespace re-ordering (alphabetized) with the
# exact same attribute values -- e.g. <w:tcW w:w="3070" w:type="dxa"/> vs
# <w:tcW w:type="dxa" w:w="3070"/>. No text, formatting, structure, or
# property value actually changed anywhere in the document, headers,
# footnotes, or styles. So the faithful edit here is a no-op on content:
# simply touch/resave the document without altering anything.
$d = $word.ActiveDocument
$d.Save()
